$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Force text storage for numeric-looking strings so Excel
    # does not auto-convert them to numbers (which would drop
    # formatting such as trailing zeros).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "62.104.22"
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("D3").Value = "2.415.84"
$ws.Range("E3").Value = "  +3.94%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "559.64"
$ws.Range("E5").Value = "  +2.69%  "
Set-TextValue "D6" "138.60"
$ws.Range("E6").Value = "  +5.92%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.413.93"
$ws.Range("E9").Value = "  +3.98%  "
$ws.Range("E10").Value = "  +3.04%  "
Set-TextValue "D11" "5.75"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("E12").Value = "  -0.16%  "
Set-TextValue "D13" "0.347"
$ws.Range("E13").Value = "  +3.47%  "
Set-TextValue "D14" "25.82"
$ws.Range("E14").Value = "  +9.38%  "
$ws.Range("D15").Value = "2.845.67"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").Value = "62.027.59"
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("D18").Value = "2.407.33"
$ws.Range("E18").Value = "  +3.82%  "
Set-TextValue "D19" "11.08"
$ws.Range("E19").Value = "  +4.72%  "
Set-TextValue "D20" "343.14"
$ws.Range("E20").Value = "  +9.28%  "
Set-TextValue "D21" "4.23"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("E23").Value = "  +0.13%  "
Set-TextValue "D24" "64.96"
$ws.Range("E24").Value = "  +1.82%  "
Set-TextValue "D25" "0.173"
$ws.Range("E26").Value = "  +0.16%  "
Set-TextValue "D27" "8.29"
$ws.Range("E27").Value = "  +5.51%  "
Set-TextValue "D28" "1.50"
$ws.Range("E28").Value = "  +11.20%  "
Set-TextValue "D29" "1.37"
$ws.Range("E29").Value = "  +14.92%  "
$ws.Range("D30").Value = "0.0₃0780"
$ws.Range("E30").Value = "  +6.62%  "
$ws.Range("E31").Value = "  +3.06%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "171.52"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D33" "6.33"
$ws.Range("E33").Value = "  +6.58%  "
$ws.Range("E34").Value = "  +3.00%  "
Set-TextValue "D35" "0.396"
$ws.Range("E35").Value = "  +3.85%  "
Set-TextValue "D36" "376.22"
$ws.Range("E36").Value = "  +16.49%  "
Set-TextValue "D37" "18.50"
$ws.Range("E37").Value = "  +3.86%  "
Set-TextValue "D38" "4.49"
$ws.Range("E38").Value = "  +10.38%  "
Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  -0.08%  "
Set-TextValue "D41" "1.67"
$ws.Range("E41").Value = "  +9.05%  "
Set-TextValue "D42" "39.07"
$ws.Range("E42").Value = "  +2.94%  "
Set-TextValue "D43" "145.42"
$ws.Range("E43").Value = "  +5.42%  "
$ws.Range("E44").Value = "  +4.80%  "
Set-TextValue "D45" "20.64"
$ws.Range("E45").Value = "  +8.24%  "
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D47" "0.587"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D48" "0.0516"
$ws.Range("E48").Value = "  +4.41%  "
Set-TextValue "D49" "17.98"
$ws.Range("E49").Value = "  +6.18%  "
Set-TextValue "D50" "0.0220"
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("D51").Value = "0.0₆0224"
$ws.Range("E51").Value = "  +4.89%  "
